# Library management system: "return/remove book" action on row 2 (1984).
# Clears the active loan (Status/Start Date/End Date/Issued By) and makes the
# book Available again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the yyyy-mm-dd date format to the Start/End date cells first (this is
# the format the "clear loan" routine always (re)applies to date cells before
# writing 0 into them).
$ws.Range("F2:G2").NumberFormat = "yyyy-mm-dd"

# ISBN cell gets re-written with the sheet's normal (non-header) style instead
# of the legacy wrapped/Times-New-Roman formatting it had before.
$ws.Range("B2").WrapText = $false
$ws.Range("B2").Font.Name = $ws.Range("A2").Font.Name

$ws.Range("A2").Value = "1984"
$ws.Range("E2").Value = "Available"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
